# Aggiornato file ore progetto al 22/01
# - new timesheet rows for 22/01/2018 (serial 43122): Mirko + Giovanni
# - small styled marker cell at F4 (mirrors the existing one at G5)
# - selection moved to F4 (no more frozen topLeftCell scroll position)
# GioH/MirkoH/Totale ore (G2/H2/J2) recalculate automatically from the
# SUMIFS formulas already on the sheet once the new D-column hours are in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Empty "marker" cell at F4, same underline-only style already used at G5.
$ws.Range("F4").Font.Underline = $true

# Row 32 - Mirko, "sistemato display con clear back scroll"
$ws.Range("A32").Value2 = 43122
$ws.Range("B32").Value2 = "Mirko"
$ws.Range("C32").Value2 = "sistemato display con clear back scroll"
$ws.Range("D32").Value2 = 0.125
$ws.Rows.Item(32).RowHeight = 28.8

# Row 33 - Giovanni, "sistemata logica calcolatrice"
$ws.Range("A33").Value2 = 43122
$ws.Range("B33").Value2 = "Giovanni"
$ws.Range("C33").Value2 = "sistemata logica calcolatrice"
$ws.Range("D33").Value2 = 0.125
$ws.Rows.Item(33).RowHeight = 28.8

# Cursor ends up on F4 (matches the saved selection in the sheet view).
$ws.Range("F4").Select() | Out-Null
